# Add a new slide (9th) at the end of the deck with a drug-ranking table,
# matching the "DiSCoVER: top drugs (cerebellar stem cell control)" slide.

$p = $ppt.ActivePresentation

# --- New slide, blank layout (matches slide8's slideLayout7 reference) ---
$slideIndex = $p.Slides.Count + 1
$newSlide = $p.Slides.Add($slideIndex, 12)

# --- Title textbox ------------------------------------------------------
$tb = $newSlide.Shapes.AddTextbox(1, 0, 0, 720, 54)
$tb.Name = "TextBox 1"
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = $false
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "DiSCoVER: top drugs (cerebellar stem cell control)"
$tb.TextFrame.TextRange.InsertBefore("`r")
$titlePara = $tb.TextFrame.TextRange.Paragraphs(2, 1)
$titlePara.Font.Size = 26
$tb.Width = 720
$tb.Height = 54

# --- Drug table ----------------------------------------------------------
$tableData = @(
    @("Drug", "Score", "Evidence", "Mechanism of action"),
    @("tl-2-105", "0.61", "+..", "Not Clinically Relevant"),
    @("gsk1070916", "0.6", "+..", "Not Clinically Relevant"),
    @("sb52334", "0.58", "+..", "Not Clinically Relevant"),
    @("gw-2580", "0.57", "+..", "Not Clinically Relevant"),
    @("gsk429286a", "0.55", "+..", "Not Clinically Relevant"),
    @("linsitinib", "0.54", "++.", "IGF-1R inhibitor"),
    @("tubastatin a", "0.53", "++.", "Not Clinically Relevant"),
    @("vx-702", "0.52", "+..", "Not Clinically Relevant"),
    @("bx-912", "0.51", "+..", "Not Clinically Relevant"),
    @("rucaparib", "0.5", "+..", "PARP inhibitor, inhibits DNA repair"),
    @("gsk319347a", "0.47", "+..", "Not Clinically Relevant"),
    @("nsc-87877", "0.46", "+..", "Not Clinically Relevant"),
    @("navitoclax", "0.46", "++.", "Bcl-2 family inhibitor: esp Bcl-xL, Bcl-2 and Bcl-w"),
    @("axitinib", "0.46", "++.", "VEGFR, c-KIT and PDGFR inhibitor"),
    @("amuvatinib", "0.46", "+..", "Not Clinically Relevant"),
    @("xmd13-2", "0.44", "+..", "Not Clinically Relevant"),
    @("staurosporine", "0.43", ".+.", "Not Clinically Relevant"),
    @("avrainvillamide", "0.42", ".+.", "Not Clinically Relevant"),
    @("talazoparib", "0.42", "+..", "Not Clinically Relevant"),
    @("bms-195614", "0.42", ".+.", "Not Clinically Relevant")
)

$numRows = $tableData.Length
$numCols = 4

$tbl = $newSlide.Shapes.AddTable($numRows, $numCols, 32.4, 61.2, 651.6, 324)
$tbl.Name = "Table 2"
$table = $tbl.Table

# Column widths (points; EMU / 12700)
$table.Columns(1).Width = 79.2
$table.Columns(2).Width = 57.6
$table.Columns(3).Width = 82.8
$table.Columns(4).Width = 432.0

# Row heights (points; EMU / 12700) - last row is very slightly taller
for ($r = 1; $r -le $numRows; $r++) {
    if ($r -eq $numRows) {
        $table.Rows($r).Height = 15.42992125984252
    } else {
        $table.Rows($r).Height = 15.428503937007873
    }
}

for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $table.Cell($r, $c)
        $cell.Shape.TextFrame.TextRange.Text = $tableData[$r - 1][$c - 1]
        if ($r -gt 1) {
            $cell.Shape.TextFrame.TextRange.Font.Size = 10.5
        }
    }
}

# Re-assert the table frame position/size after population (row
# auto-fit during text entry can perturb the authored extent).
$tbl.Left = 32.4
$tbl.Top = 61.2
$tbl.Width = 651.6
$tbl.Height = 324
